$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 367, pushing the existing rows 367-402
# down to 369-404 (dimension grows from A1:R402 to A1:R404).
$ws.Range("367:368").Insert()

# New row 367: Betarraga / Primera, week of 2022-07-27 (serial 44769)
$ws.Range("A367").Value = 7
$ws.Range("B367").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C367").Value = "Ñuble"
$ws.Range("D367").Value = 44769
$ws.Range("E367").Value = 16
$ws.Range("F367").Value = 100114014
$ws.Range("G367").Value = "Betarraga"
$ws.Range("H367").Value = "Sin especificar"
$ws.Range("I367").Value = "Primera"
$ws.Range("J367").Value = 200
$ws.Range("K367").Value = 700
$ws.Range("L367").Value = 800
$ws.Range("M367").Value = 750
$ws.Range("N367").Value = "$/paquete 5 unidades"
$ws.Range("O367").Value = "Provincia de Diguillín"
$ws.Range("P367").Value = 150
$ws.Range("Q367").Value = 5
$ws.Range("R367").Value = "Hortaliza"

# New row 368: Betarraga / Segunda, same week (serial 44769)
$ws.Range("A368").Value = 7
$ws.Range("B368").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C368").Value = "Ñuble"
$ws.Range("D368").Value = 44769
$ws.Range("E368").Value = 16
$ws.Range("F368").Value = 100114014
$ws.Range("G368").Value = "Betarraga"
$ws.Range("H368").Value = "Sin especificar"
$ws.Range("I368").Value = "Segunda"
$ws.Range("J368").Value = 150
$ws.Range("K368").Value = 600
$ws.Range("L368").Value = 600
$ws.Range("M368").Value = 600
$ws.Range("N368").Value = "$/paquete 5 unidades"
$ws.Range("O368").Value = "Provincia de Diguillín"
$ws.Range("P368").Value = 120
$ws.Range("Q368").Value = 5
$ws.Range("R368").Value = "Hortaliza"
